# Add a "Slovakia" market tab to the workbook, built from a copy of the
# "Portugal" sheet (same layout/repeater list minus the CH-only models),
# with its own User Story reference.

$wb = $excel.ActiveWorkbook
$portugal = $wb.Worksheets.Item("Portugal")

# Portugal is no longer the active tab once Slovakia exists; its selection
# reverts to the whole sheet (equivalent to Ctrl+A) instead of B2.
$null = $portugal.Cells.Select()

# Duplicate "Portugal" -> new sheet lands right after it and becomes active.
$null = $portugal.Copy($null, $portugal)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Slovakia"

# Slovakia's repeater list doesn't include P32AR / P32DR - drop those two
# rows (bottom-up so the second delete doesn't need an index shift).
$null = $ws.Rows.Item(17).Delete()
$null = $ws.Rows.Item(16).Delete()

# Slovakia-specific User Story reference, in the default (unstyled) cell.
$ws.Range("B4").Style = "Normal"
$ws.Range("B4").Value = "NGC-2930/T3176/T3179"

# Re-fit the rows that used to be sized for the old wrapped content.
$null = $ws.Rows.Item(3).AutoFit()
$null = $ws.Rows.Item(4).AutoFit()
$null = $ws.Rows.Item(5).AutoFit()

$null = $ws.Range("B4").Select()
